# Weekly update: insert two new daily price rows at the top of the
# Plátano logged data (rows 629-630) and push the existing history
# down by two rows (old 629..728 -> new 631..730).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 629..728 down by two rows, preserving formatting
# (in particular the date number format on column D).
$ws.Rows.Item(629).Resize(2).Insert()

# New row 629
$row629 = @(4, "Feria Lagunitas de Puerto Montt", "Los Lagos", 44951, 10, "Fruta", 100108, "Tropicales y subtropicales", 100108006, "Plátano", "Sin especificar", "Pintón", 300, 25000, 25000, 25000, "$/caja 20 kilos", "Ecuador", 1250, 20)
for ($i = 0; $i -lt $row629.Length; $i++) {
    $ws.Cells.Item(629, $i + 1).Value = $row629[$i]
}

# New row 630
$row630 = @(4, "Feria Lagunitas de Puerto Montt", "Los Lagos", 44951, 10, "Fruta", 100108, "Tropicales y subtropicales", 100108006, "Plátano", "Sin especificar", "Primera Pintón", 600, 26000, 27000, 26500, "$/caja 20 kilos", "Ecuador", 1325, 20)
for ($i = 0; $i -lt $row630.Length; $i++) {
    $ws.Cells.Item(630, $i + 1).Value = $row630[$i]
}
